# Updates Price (col D) and Volume(1h) (col E) values for the cryptos sheet
# based on the latest scrape (GitHub Actions commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.533.10'
$ws.Range('E2').Value = '  +0.82%  '
$ws.Range('D3').Value = '1.879.61'
$ws.Range('E3').Value = '  +1.16%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7186'
$ws.Range('E5').Value = '  +2.60%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '242.08'
$ws.Range('E6').Value = '  +1.68%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.07944'
$ws.Range('E8').Value = '  +1.06%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3103'
$ws.Range('E9').Value = '  +2.56%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '25.44'
$ws.Range('E10').Value = '  +4.14%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.08277'
$ws.Range('E11').Value = '  +1.35%  '
$ws.Range('D12').Value = '1.893.57'
$ws.Range('E12').Value = '  +1.43%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.7313'
$ws.Range('E13').Value = '  +3.56%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.286'
$ws.Range('E14').Value = '  +1.50%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '91.32'
$ws.Range('E15').Value = '  +2.05%  '
$ws.Range('D16').Value = '29.531.47'
$ws.Range('E16').Value = '  +0.62%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '5.910'
$ws.Range('E17').Value = '  +2.02%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '246.34'
$ws.Range('E18').Value = '  +3.92%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000007872'
$ws.Range('E19').Value = '  +0.79%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.35'
$ws.Range('E20').Value = '  +1.15%  '
$ws.Range('D21').Value = '2.125.48'
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '8.073'
$ws.Range('E22').Value = '  +6.78%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.002'
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.001'
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.1624'
$ws.Range('E25').Value = '  +14.77%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '163.53'
$ws.Range('E26').Value = '  +0.81%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.062'
$ws.Range('E27').Value = '  +1.95%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '18.34'
$ws.Range('E28').Value = '  +1.49%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.353'
$ws.Range('E29').Value = '  -3.26%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.496'
$ws.Range('E30').Value = '  +1.18%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.403'
$ws.Range('E31').Value = '  +2.47%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.114'
$ws.Range('E32').Value = '  +2.08%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05227'
$ws.Range('E33').Value = '  +0.85%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.952'
$ws.Range('E34').Value = '  +2.49%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.200'
$ws.Range('E35').Value = '  +1.85%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7292'
$ws.Range('E36').Value = '  +2.79%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.678'
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01874'
$ws.Range('E38').Value = '  +1.40%  '
$ws.Range('D39').Value = '1.205.52'
$ws.Range('E39').Value = '  +5.69%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.706'
$ws.Range('E40').Value = '  +0.80%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.9113'
$ws.Range('E41').Value = '  -1.19%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.181'
$ws.Range('E42').Value = '  +3.97%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '73.64'
$ws.Range('E43').Value = '  +4.74%  '
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '102.46'
$ws.Range('E45').Value = '  -0.33%  '
$ws.Range('D46').Value = '2.022.15'
$ws.Range('E46').Value = '  -0.24%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5294'
$ws.Range('E47').Value = '  -0.46%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.803'
$ws.Range('E48').Value = '  +3.43%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.936'
$ws.Range('E49').Value = '  +9.68%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '9.317'
$ws.Range('E50').Value = '  +1.39%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.4332'
$ws.Range('E51').Value = '  +2.19%  '
